$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B4 changes from "Test123" to "Test"
$ws.Range("B4").Value = "Test"

# Update the active selection to L10 (matches the saved view state in the diff)
$ws.Range("L10").Select()
